# Update transition-matrix probabilities on Sheet1 (Sacramento St._A) to
# reflect the recalculated values from games pulled March 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.2118055555555556
$ws.Cells.Item(2, 3).Value = 0.5173611111111112
$ws.Cells.Item(2, 10).Value = 0.02777777777777778
$ws.Cells.Item(2, 16).Value = 0.1597222222222222
$ws.Cells.Item(2, 19).Value = 0.08333333333333333
$ws.Cells.Item(3, 3).Value = 0.0379746835443038
$ws.Cells.Item(3, 10).Value = 0.01265822784810127
$ws.Cells.Item(3, 16).Value = 0.740506329113924
$ws.Cells.Item(3, 19).Value = 0.2088607594936709
$ws.Cells.Item(4, 10).Value = 0.02222222222222222
$ws.Cells.Item(4, 15).Value = 0.02222222222222222
$ws.Cells.Item(4, 16).Value = 0.6666666666666666
$ws.Cells.Item(4, 19).Value = 0.2888888888888889
$ws.Cells.Item(6, 2).Value = 0.05405405405405406
$ws.Cells.Item(6, 4).Value = 0.01158301158301158
$ws.Cells.Item(6, 6).Value = 0.06563706563706563
$ws.Cells.Item(6, 10).Value = 0.2432432432432433
$ws.Cells.Item(6, 15).Value = 0.03474903474903475
$ws.Cells.Item(6, 17).Value = 0.1467181467181467
$ws.Cells.Item(6, 18).Value = 0.06563706563706563
$ws.Cells.Item(6, 19).Value = 0.3783783783783784
$ws.Cells.Item(7, 2).Value = 0.09944751381215469
$ws.Cells.Item(7, 4).Value = 0.02209944751381215
$ws.Cells.Item(7, 6).Value = 0.09944751381215469
$ws.Cells.Item(7, 10).Value = 0.1491712707182321
$ws.Cells.Item(7, 15).Value = 0.03314917127071823
$ws.Cells.Item(7, 17).Value = 0.1104972375690608
$ws.Cells.Item(7, 18).Value = 0.09392265193370165
$ws.Cells.Item(7, 19).Value = 0.3922651933701657
$ws.Cells.Item(8, 2).Value = 0.07468879668049792
$ws.Cells.Item(8, 4).Value = 0.01867219917012448
$ws.Cells.Item(8, 6).Value = 0.07261410788381743
$ws.Cells.Item(8, 10).Value = 0.1535269709543569
$ws.Cells.Item(8, 15).Value = 0.03112033195020747
$ws.Cells.Item(8, 17).Value = 0.1701244813278008
$ws.Cells.Item(8, 18).Value = 0.09336099585062241
$ws.Cells.Item(8, 19).Value = 0.3858921161825726
$ws.Cells.Item(9, 2).Value = 0.08878504672897196
$ws.Cells.Item(9, 4).Value = 0.009345794392523364
$ws.Cells.Item(9, 6).Value = 0.07476635514018691
$ws.Cells.Item(9, 10).Value = 0.1308411214953271
$ws.Cells.Item(9, 15).Value = 0.04205607476635514
$ws.Cells.Item(9, 17).Value = 0.1308411214953271
$ws.Cells.Item(9, 18).Value = 0.09813084112149532
$ws.Cells.Item(9, 19).Value = 0.4252336448598131
$ws.Cells.Item(10, 2).Value = 0.1018587360594796
$ws.Cells.Item(10, 4).Value = 0.02156133828996282
$ws.Cells.Item(10, 5).Value = 0.001486988847583643
$ws.Cells.Item(10, 6).Value = 0.07657992565055761
$ws.Cells.Item(10, 10).Value = 0.1353159851301115
$ws.Cells.Item(10, 15).Value = 0.02304832713754647
$ws.Cells.Item(10, 17).Value = 0.1955390334572491
$ws.Cells.Item(10, 18).Value = 0.08847583643122676
$ws.Cells.Item(10, 19).Value = 0.3561338289962825
$ws.Cells.Item(11, 7).Value = 0.1160409556313993
$ws.Cells.Item(11, 10).Value = 0.10580204778157
$ws.Cells.Item(11, 11).Value = 0.1877133105802048
$ws.Cells.Item(11, 12).Value = 0.5767918088737202
$ws.Cells.Item(11, 19).Value = 0.0136518771331058
$ws.Cells.Item(12, 7).Value = 0.7352941176470589
$ws.Cells.Item(12, 10).Value = 0.2058823529411765
$ws.Cells.Item(12, 11).Value = 0.005882352941176471
$ws.Cells.Item(12, 12).Value = 0.02352941176470588
$ws.Cells.Item(12, 19).Value = 0.02941176470588235
$ws.Cells.Item(13, 7).Value = 0.6170212765957447
$ws.Cells.Item(13, 10).Value = 0.2765957446808511
$ws.Cells.Item(13, 19).Value = 0.1063829787234043
$ws.Cells.Item(15, 6).Value = 0.01515151515151515
$ws.Cells.Item(15, 8).Value = 0.1553030303030303
$ws.Cells.Item(15, 9).Value = 0.07575757575757576
$ws.Cells.Item(15, 10).Value = 0.3181818181818182
$ws.Cells.Item(15, 11).Value = 0.05681818181818182
$ws.Cells.Item(15, 13).Value = 0.01515151515151515
$ws.Cells.Item(15, 15).Value = 0.05303030303030303
$ws.Cells.Item(15, 19).Value = 0.3106060606060606
$ws.Cells.Item(16, 6).Value = 0.02127659574468085
$ws.Cells.Item(16, 8).Value = 0.2180851063829787
$ws.Cells.Item(16, 9).Value = 0.09042553191489362
$ws.Cells.Item(16, 10).Value = 0.4095744680851064
$ws.Cells.Item(16, 11).Value = 0.06914893617021277
$ws.Cells.Item(16, 13).Value = 0.02659574468085106
$ws.Cells.Item(16, 15).Value = 0.04787234042553191
$ws.Cells.Item(16, 19).Value = 0.1170212765957447
$ws.Cells.Item(17, 6).Value = 0.01635514018691589
$ws.Cells.Item(17, 8).Value = 0.1869158878504673
$ws.Cells.Item(17, 9).Value = 0.1121495327102804
$ws.Cells.Item(17, 10).Value = 0.4042056074766355
$ws.Cells.Item(17, 11).Value = 0.09579439252336448
$ws.Cells.Item(17, 13).Value = 0.02336448598130841
$ws.Cells.Item(17, 15).Value = 0.07476635514018691
$ws.Cells.Item(17, 19).Value = 0.08644859813084112
$ws.Cells.Item(18, 6).Value = 0.0273972602739726
$ws.Cells.Item(18, 8).Value = 0.2191780821917808
$ws.Cells.Item(18, 9).Value = 0.0821917808219178
$ws.Cells.Item(18, 10).Value = 0.4611872146118721
$ws.Cells.Item(18, 11).Value = 0.0821917808219178
$ws.Cells.Item(18, 13).Value = 0.0091324200913242
$ws.Cells.Item(18, 15).Value = 0.0502283105022831
$ws.Cells.Item(18, 19).Value = 0.0684931506849315
$ws.Cells.Item(19, 6).Value = 0.01648351648351648
$ws.Cells.Item(19, 8).Value = 0.2182103610675039
$ws.Cells.Item(19, 9).Value = 0.08869701726844584
$ws.Cells.Item(19, 10).Value = 0.3618524332810047
$ws.Cells.Item(19, 11).Value = 0.1185243328100471
$ws.Cells.Item(19, 13).Value = 0.02197802197802198
$ws.Cells.Item(19, 15).Value = 0.07378335949764521
$ws.Cells.Item(19, 19).Value = 0.1004709576138148
